# "mostly contact angle stuff"
# Add the standard-deviation summary row under the existing AVERAGE row,
# mirroring the new B16 = STDEV(B9:B13) cell added to the worksheet, and
# move the active selection to B9 (where the data entry now focuses).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Formula = "=STDEV(B9:B13)"

$ws.Range("B9").Select() | Out-Null
